# Apply a cyclic rotation of the weekly price data across rows 2, 3 and 4:
#   new Row2 (D,J,K,L,M,P) = old Row4
#   new Row3 (D,J,K,L,M,P) = old Row2
#   new Row4 (D,J,K,L,M,P) = old Row3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "P")

# Capture the original (before-edit) values for the affected columns/rows.
$orig = @{}
foreach ($col in $cols) {
    $orig[$col] = @{
        2 = $ws.Range("${col}2").Value2
        3 = $ws.Range("${col}3").Value2
        4 = $ws.Range("${col}4").Value2
    }
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $orig[$col][4]
    $ws.Range("${col}3").Value2 = $orig[$col][2]
    $ws.Range("${col}4").Value2 = $orig[$col][3]
}
